$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "HK_G_acc_LG"

$values = @(
    76.511954992967645,
    76.511954992967645,
    76.230661040787624,
    76.511954992967645,
    76.371308016877634,
    76.371308016877634,
    76.371308016877634,
    76.090014064697613,
    76.371308016877634,
    76.371308016877634,
    75.949367088607602,
    76.511954992967645,
    76.511954992967645,
    76.230661040787624,
    76.511954992967645,
    76.511954992967645,
    76.511954992967645,
    76.511954992967645,
    76.230661040787624,
    76.230661040787624,
    76.230661040787624,
    76.090014064697613,
    76.230661040787624,
    75.808720112517577,
    77.355836849507725,
    76.511954992967645,
    77.355836849507725,
    76.090014064697613,
    76.230661040787624,
    76.511954992967645,
    76.371308016877634,
    76.511954992967645,
    76.511954992967645,
    76.65260196905767,
    75.668073136427566,
    75.386779184247544,
    76.090014064697613,
    77.074542897327703,
    77.215189873417728,
    76.511954992967645,
    76.371308016877634,
    76.371308016877634,
    76.371308016877634,
    76.511954992967645,
    76.371308016877634,
    76.511954992967645,
    76.371308016877634,
    76.511954992967645
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
